# Update inventory and BOMs
# - Bump the "Cnt" column on the BGB741L7ESD and TQP3M9037 amp BOM sheets
#   to reflect updated inventory counts.
# - Append a new revision-history row to each of those two sheets noting the
#   latest assembled-board change, dated 2019-11-13 (serial 43782).

$wb  = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item("amp-LNA-BGB741L7ESD")
$ws3 = $wb.Worksheets.Item("amp-LNA-TQP3M9037")

# --- amp-LNA-BGB741L7ESD (sheet2) Cnt column updates ---
$ws2.Range("A8").Value  = 3
$ws2.Range("A9").Value  = 4
$ws2.Range("A11").Value = 2
$ws2.Range("A12").Value = 4
$ws2.Range("A13").Value = 4
$ws2.Range("A14").Value = 2
$ws2.Range("A15").Value = 2
$ws2.Range("A18").Value = 3
$ws2.Range("A19").Value = 2
$ws2.Range("A20").Value = 2
$ws2.Range("A21").Value = 4

# --- amp-LNA-TQP3M9037 (sheet3) Cnt column updates ---
$ws3.Range("A8").Value  = 3
$ws3.Range("A9").Value  = 2
$ws3.Range("A10").Value = 2
$ws3.Range("A12").Value = 2
$ws3.Range("A13").Value = 2
$ws3.Range("A15").Value = 2
$ws3.Range("A16").Value = 2
$ws3.Range("A17").Value = 2
$ws3.Range("A19").Value = 2
$ws3.Range("A20").Value = 2
$ws3.Range("A21").Value = 4

# --- New revision-history rows ---
# Write the TQP3M9037 note first so its shared string lands at the lower
# index (161), matching BGB741L7ESD's longer note at index 162.
$ws3.Range("B32").Value = 43782
$ws3.Range("C32").Value = "Assembled 1x board, NP JP2"

$ws2.Range("B34").Value = 43782
$ws2.Range("C34").Value = "Assembled 1x board, 0603 jumper for Q1, NP R3, NP JP2"

# Match the existing yyyy-mm-dd date formatting used by the rest of the
# Revision History table (style index 8) by copying format-only from the
# row above on each sheet.
$ws3.Range("B31").Copy()
$ws3.Range("B32").PasteSpecial(-4122)

$ws2.Range("B33").Copy()
$ws2.Range("B34").PasteSpecial(-4122)
